$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("J2").Value = 3.75
$ws.Range("K2").Value = 3.85
$ws.Range("X2").Value = 12
$ws.Range("Z2").Value = 36

# Row 3
$ws.Range("P3").Value = 1.76
$ws.Range("Q3").Value = 2.28

# Row 4
$ws.Range("F4").Value = 5.7
$ws.Range("G4").Value = 5.8
$ws.Range("H4").Value = 1.83
$ws.Range("I4").Value = 1.85
$ws.Range("J4").Value = 3.55
$ws.Range("L4").Value = 1.5
$ws.Range("N4").Value = 2.98
$ws.Range("O4").Value = 1.48
$ws.Range("Q4").Value = 2.44
$ws.Range("S4").Value = 4.9
$ws.Range("U4").Value = 1.73
$ws.Range("V4").Value = 2.18
$ws.Range("W4").Value = 1.2
$ws.Range("X4").Value = 9.800000000000001
$ws.Range("Z4").Value = 9.199999999999999
$ws.Range("AA4").Value = 19
$ws.Range("AB4").Value = 15
$ws.Range("AC4").Value = 8
$ws.Range("AD4").Value = 10.5
$ws.Range("AE4").Value = 24
$ws.Range("AF4").Value = 40
$ws.Range("AG4").Value = 23
$ws.Range("AH4").Value = 27
$ws.Range("AI4").Value = 55
$ws.Range("AJ4").Value = 180
$ws.Range("AK4").Value = 100
$ws.Range("AL4").Value = 120
$ws.Range("AM4").Value = 200
$ws.Range("AN4").Value = 170
$ws.Range("AO4").Value = 17

# Row 5
$ws.Range("K5").Value = 3.95
$ws.Range("Q5").Value = 2.04
$ws.Range("AB5").Value = 21
$ws.Range("AC5").Value = 9.4

# Row 6
$ws.Range("F6").Value = 2.62
$ws.Range("G6").Value = 2.64
$ws.Range("H6").Value = 2.9
$ws.Range("I6").Value = 2.94
$ws.Range("J6").Value = 3.6
$ws.Range("P6").Value = 1.92
$ws.Range("Q6").Value = 1.9

# Row 7
$ws.Range("F7").Value = 3.25
$ws.Range("I7").Value = 2.42
$ws.Range("N7").Value = 3.85
$ws.Range("P7").Value = 1.94
$ws.Range("Q7").Value = 2.04
$ws.Range("T7").Value = 1.78
$ws.Range("AA7").Value = 34
$ws.Range("AF7").Value = 23
$ws.Range("AH7").Value = 16.5
$ws.Range("AI7").Value = 42
$ws.Range("AJ7").Value = 60
$ws.Range("AK7").Value = 36
$ws.Range("AM7").Value = 80

# Row 8
$ws.Range("Z8").Value = 55
$ws.Range("AG8").Value = 13
